# Weekly update: insert one new price record as row 537, pushing the
# existing rows 537-597 down to 538-598 (dimension grows from R597 to R598).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 537; everything below shifts down one.
$ws.Rows(537).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(537, 1).Value = 6
$ws.Cells.Item(537, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(537, 3).Value = "Metropolitana"
$ws.Cells.Item(537, 4).Value = 45124
$ws.Cells.Item(537, 5).Value = 13
$ws.Cells.Item(537, 6).Value = 100112032
$ws.Cells.Item(537, 7).Value = "Zapallo italiano"
$ws.Cells.Item(537, 8).Value = "Sin especificar"
$ws.Cells.Item(537, 9).Value = "Primera"
$ws.Cells.Item(537, 10).Value = 290
$ws.Cells.Item(537, 11).Value = 15000
$ws.Cells.Item(537, 12).Value = 16000
$ws.Cells.Item(537, 13).Value = 15586
$ws.Cells.Item(537, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(537, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(537, 16).Value = 312
$ws.Cells.Item(537, 17).Value = 50
$ws.Cells.Item(537, 18).Value = "Hortaliza"
